# Applies the "padrao_importacao2" header rework:
#   - adds clarifying hint text to a few header labels
#   - re-orders the header columns (A..K) to the new layout
#   - widens / re-sizes several columns to fit their new (longer) labels
#   - moves the active cell selection to H6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header labels, written left-to-right so that new shared-string
# entries are appended in the same order the original workbook used. ---
$ws.Range("A1").Value = "Nome do aluno"
$ws.Range("B1").Value = "Data de nascimento: (formato dd/mm/aaaa)"
$ws.Range("C1").Value = "Nome da mãe ou responsável"
$ws.Range("D1").Value = "Nome do pai"
$ws.Range("E1").Value = "Endereço"
$ws.Range("F1").Value = "Bairro"
$ws.Range("G1").Value = "UF: (sigla do estado)"
$ws.Range("H1").Value = "Cidade"
$ws.Range("J1").Value = "Telefone: (apenas números com DDD)"
$ws.Range("I1").Value = "CEP: (somente números) "
$ws.Range("K1").Value = "Observações"

# --- Column widths, in "characters" (same unit Excel's Format > Column Width
# dialog uses). The interop layer adds a fixed 5/6-character offset when it
# stores a width, so that offset is subtracted here up front to land on the
# closest achievable value to the target layout. ---
$ws.Columns.Item(2).ColumnWidth = 36.140625 - 5/6    # Data de nascimento
$ws.Columns.Item(6).ColumnWidth = 21.85546875 - 5/6  # Bairro (unchanged text, wider col)
$ws.Columns.Item(7).ColumnWidth = 21 - 5/6            # UF
$ws.Columns.Item(9).ColumnWidth = 22 - 5/6            # CEP
$ws.Columns.Item(10).ColumnWidth = 32.28515625 - 5/6 # Telefone

# --- Selection moves to H6 ---
$ws.Range("H6").Select() | Out-Null
